# Update quarterly report: roll the quarter columns forward by one
# (drop Q2 1400/06, shift the rest left, append Q3 1401/09) and refresh
# the copyright year + the numeric figures for the newly-shifted quarters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copyright footer -------------------------------------------------
$ws.Range("B3").Value = "Copyright @2015 - 2023"

# --- Quarter header labels (both tables) -------------------------------
$headers = @("فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09")
$ws.Range("E8").Value  = $headers[0]
$ws.Range("F8").Value  = $headers[1]
$ws.Range("G8").Value  = $headers[2]
$ws.Range("H8").Value  = $headers[3]
$ws.Range("I8").Value  = $headers[4]

$ws.Range("E24").Value = $headers[0]
$ws.Range("F24").Value = $headers[1]
$ws.Range("G24").Value = $headers[2]
$ws.Range("H24").Value = $headers[3]
$ws.Range("I24").Value = $headers[4]

# --- Numeric rows: shift values left one quarter, add the new quarter --
$rowData = @{
    10 = @(59868, 45202, 46067, 60929, 66233)
    12 = @(0, 0, 0, 0, 182377)
    14 = @(1261, 1397, 1311, 1033, 2559)
    15 = @(242, -42, 95, 134, 132)
    16 = @(2351, 2327, 2374, 3466, 3471)
    17 = @(28575, 29145, 39914, 39471, 40677)
    19 = @(11492, 7415, 9238, 21048, 12237)
    20 = @(103789, 85444, 98999, 126081, 307686)
    26 = @(259, 256, 262, 261, 260)
    27 = @(479, 468, 508, 509, 505)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    $ws.Cells.Item($row, 5).Value = $values[0]  # E
    $ws.Cells.Item($row, 6).Value = $values[1]  # F
    $ws.Cells.Item($row, 7).Value = $values[2]  # G
    $ws.Cells.Item($row, 8).Value = $values[3]  # H
    $ws.Cells.Item($row, 9).Value = $values[4]  # I
}
